# Auto-generated: reorder category/group code-name columns (D,E,F,G) and relabel header
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object "object[,]" 235,4

$data[0,0] = "codeforiati:group-name"; $data[0,1] = "codeforiati:category-name"; $data[0,2] = "codeforiati:group-code"; $data[0,3] = "codeforiati:category-code"
$data[1,0] = "Education"; $data[1,1] = "Education, Level Unspecified"; $data[1,2] = "110"; $data[1,3] = "111"
$data[2,0] = "Education"; $data[2,1] = "Education, Level Unspecified"; $data[2,2] = "110"; $data[2,3] = "111"
$data[3,0] = "Education"; $data[3,1] = "Education, Level Unspecified"; $data[3,2] = "110"; $data[3,3] = "111"
$data[4,0] = "Education"; $data[4,1] = "Education, Level Unspecified"; $data[4,2] = "110"; $data[4,3] = "111"
$data[5,0] = "Education"; $data[5,1] = "Basic Education"; $data[5,2] = "110"; $data[5,3] = "112"
$data[6,0] = "Education"; $data[6,1] = "Basic Education"; $data[6,2] = "110"; $data[6,3] = "112"
$data[7,0] = "Education"; $data[7,1] = "Basic Education"; $data[7,2] = "110"; $data[7,3] = "112"
$data[8,0] = "Education"; $data[8,1] = "Basic Education"; $data[8,2] = "110"; $data[8,3] = "112"
$data[9,0] = "Education"; $data[9,1] = "Basic Education"; $data[9,2] = "110"; $data[9,3] = "112"
$data[10,0] = "Education"; $data[10,1] = "Basic Education"; $data[10,2] = "110"; $data[10,3] = "112"
$data[11,0] = "Education"; $data[11,1] = "Basic Education"; $data[11,2] = "110"; $data[11,3] = "112"
$data[12,0] = "Education"; $data[12,1] = "Secondary Education"; $data[12,2] = "110"; $data[12,3] = "113"
$data[13,0] = "Education"; $data[13,1] = "Secondary Education"; $data[13,2] = "110"; $data[13,3] = "113"
$data[14,0] = "Education"; $data[14,1] = "Post-Secondary Education"; $data[14,2] = "110"; $data[14,3] = "114"
$data[15,0] = "Education"; $data[15,1] = "Post-Secondary Education"; $data[15,2] = "110"; $data[15,3] = "114"
$data[16,0] = "Health"; $data[16,1] = "Health, General"; $data[16,2] = "120"; $data[16,3] = "121"
$data[17,0] = "Health"; $data[17,1] = "Health, General"; $data[17,2] = "120"; $data[17,3] = "121"
$data[18,0] = "Health"; $data[18,1] = "Health, General"; $data[18,2] = "120"; $data[18,3] = "121"
$data[19,0] = "Health"; $data[19,1] = "Health, General"; $data[19,2] = "120"; $data[19,3] = "121"
$data[20,0] = "Health"; $data[20,1] = "Basic Health"; $data[20,2] = "120"; $data[20,3] = "122"
$data[21,0] = "Health"; $data[21,1] = "Basic Health"; $data[21,2] = "120"; $data[21,3] = "122"
$data[22,0] = "Health"; $data[22,1] = "Basic Health"; $data[22,2] = "120"; $data[22,3] = "122"
$data[23,0] = "Health"; $data[23,1] = "Basic Health"; $data[23,2] = "120"; $data[23,3] = "122"
$data[24,0] = "Health"; $data[24,1] = "Basic Health"; $data[24,2] = "120"; $data[24,3] = "122"
$data[25,0] = "Health"; $data[25,1] = "Basic Health"; $data[25,2] = "120"; $data[25,3] = "122"
$data[26,0] = "Health"; $data[26,1] = "Basic Health"; $data[26,2] = "120"; $data[26,3] = "122"
$data[27,0] = "Health"; $data[27,1] = "Basic Health"; $data[27,2] = "120"; $data[27,3] = "122"
$data[28,0] = "Health"; $data[28,1] = "Basic Health"; $data[28,2] = "120"; $data[28,3] = "122"
$data[29,0] = "Health"; $data[29,1] = "Non-communicable diseases (NCDs)"; $data[29,2] = "120"; $data[29,3] = "123"
$data[30,0] = "Health"; $data[30,1] = "Non-communicable diseases (NCDs)"; $data[30,2] = "120"; $data[30,3] = "123"
$data[31,0] = "Health"; $data[31,1] = "Non-communicable diseases (NCDs)"; $data[31,2] = "120"; $data[31,3] = "123"
$data[32,0] = "Health"; $data[32,1] = "Non-communicable diseases (NCDs)"; $data[32,2] = "120"; $data[32,3] = "123"
$data[33,0] = "Health"; $data[33,1] = "Non-communicable diseases (NCDs)"; $data[33,2] = "120"; $data[33,3] = "123"
$data[34,0] = "Health"; $data[34,1] = "Non-communicable diseases (NCDs)"; $data[34,2] = "120"; $data[34,3] = "123"
$data[35,0] = "Population Policies/Programmes & Reproductive Health"; $data[35,1] = "Population Policies/Programmes & Reproductive Health"; $data[35,2] = "130"; $data[35,3] = "130"
$data[36,0] = "Population Policies/Programmes & Reproductive Health"; $data[36,1] = "Population Policies/Programmes & Reproductive Health"; $data[36,2] = "130"; $data[36,3] = "130"
$data[37,0] = "Population Policies/Programmes & Reproductive Health"; $data[37,1] = "Population Policies/Programmes & Reproductive Health"; $data[37,2] = "130"; $data[37,3] = "130"
$data[38,0] = "Population Policies/Programmes & Reproductive Health"; $data[38,1] = "Population Policies/Programmes & Reproductive Health"; $data[38,2] = "130"; $data[38,3] = "130"
$data[39,0] = "Population Policies/Programmes & Reproductive Health"; $data[39,1] = "Population Policies/Programmes & Reproductive Health"; $data[39,2] = "130"; $data[39,3] = "130"
$data[40,0] = "Water Supply & Sanitation"; $data[40,1] = "Water Supply & Sanitation"; $data[40,2] = "140"; $data[40,3] = "140"
$data[41,0] = "Water Supply & Sanitation"; $data[41,1] = "Water Supply & Sanitation"; $data[41,2] = "140"; $data[41,3] = "140"
$data[42,0] = "Water Supply & Sanitation"; $data[42,1] = "Water Supply & Sanitation"; $data[42,2] = "140"; $data[42,3] = "140"
$data[43,0] = "Water Supply & Sanitation"; $data[43,1] = "Water Supply & Sanitation"; $data[43,2] = "140"; $data[43,3] = "140"
$data[44,0] = "Water Supply & Sanitation"; $data[44,1] = "Water Supply & Sanitation"; $data[44,2] = "140"; $data[44,3] = "140"
$data[45,0] = "Water Supply & Sanitation"; $data[45,1] = "Water Supply & Sanitation"; $data[45,2] = "140"; $data[45,3] = "140"
$data[46,0] = "Water Supply & Sanitation"; $data[46,1] = "Water Supply & Sanitation"; $data[46,2] = "140"; $data[46,3] = "140"
$data[47,0] = "Water Supply & Sanitation"; $data[47,1] = "Water Supply & Sanitation"; $data[47,2] = "140"; $data[47,3] = "140"
$data[48,0] = "Water Supply & Sanitation"; $data[48,1] = "Water Supply & Sanitation"; $data[48,2] = "140"; $data[48,3] = "140"
$data[49,0] = "Water Supply & Sanitation"; $data[49,1] = "Water Supply & Sanitation"; $data[49,2] = "140"; $data[49,3] = "140"
$data[50,0] = "Water Supply & Sanitation"; $data[50,1] = "Water Supply & Sanitation"; $data[50,2] = "140"; $data[50,3] = "140"
$data[51,0] = "Government & Civil Society"; $data[51,1] = "Government & Civil Society-general"; $data[51,2] = "150"; $data[51,3] = "151"
$data[52,0] = "Government & Civil Society"; $data[52,1] = "Government & Civil Society-general"; $data[52,2] = "150"; $data[52,3] = "151"
$data[53,0] = "Government & Civil Society"; $data[53,1] = "Government & Civil Society-general"; $data[53,2] = "150"; $data[53,3] = "151"
$data[54,0] = "Government & Civil Society"; $data[54,1] = "Government & Civil Society-general"; $data[54,2] = "150"; $data[54,3] = "151"
$data[55,0] = "Government & Civil Society"; $data[55,1] = "Government & Civil Society-general"; $data[55,2] = "150"; $data[55,3] = "151"
$data[56,0] = "Government & Civil Society"; $data[56,1] = "Government & Civil Society-general"; $data[56,2] = "150"; $data[56,3] = "151"
$data[57,0] = "Government & Civil Society"; $data[57,1] = "Government & Civil Society-general"; $data[57,2] = "150"; $data[57,3] = "151"
$data[58,0] = "Government & Civil Society"; $data[58,1] = "Government & Civil Society-general"; $data[58,2] = "150"; $data[58,3] = "151"
$data[59,0] = "Government & Civil Society"; $data[59,1] = "Government & Civil Society-general"; $data[59,2] = "150"; $data[59,3] = "151"
$data[60,0] = "Government & Civil Society"; $data[60,1] = "Government & Civil Society-general"; $data[60,2] = "150"; $data[60,3] = "151"
$data[61,0] = "Government & Civil Society"; $data[61,1] = "Government & Civil Society-general"; $data[61,2] = "150"; $data[61,3] = "151"
$data[62,0] = "Government & Civil Society"; $data[62,1] = "Government & Civil Society-general"; $data[62,2] = "150"; $data[62,3] = "151"
$data[63,0] = "Government & Civil Society"; $data[63,1] = "Government & Civil Society-general"; $data[63,2] = "150"; $data[63,3] = "151"
$data[64,0] = "Government & Civil Society"; $data[64,1] = "Government & Civil Society-general"; $data[64,2] = "150"; $data[64,3] = "151"
$data[65,0] = "Government & Civil Society"; $data[65,1] = "Government & Civil Society-general"; $data[65,2] = "150"; $data[65,3] = "151"
$data[66,0] = "Government & Civil Society"; $data[66,1] = "Government & Civil Society-general"; $data[66,2] = "150"; $data[66,3] = "151"
$data[67,0] = "Government & Civil Society"; $data[67,1] = "Conflict, Peace & Security"; $data[67,2] = "150"; $data[67,3] = "152"
$data[68,0] = "Government & Civil Society"; $data[68,1] = "Conflict, Peace & Security"; $data[68,2] = "150"; $data[68,3] = "152"
$data[69,0] = "Government & Civil Society"; $data[69,1] = "Conflict, Peace & Security"; $data[69,2] = "150"; $data[69,3] = "152"
$data[70,0] = "Government & Civil Society"; $data[70,1] = "Conflict, Peace & Security"; $data[70,2] = "150"; $data[70,3] = "152"
$data[71,0] = "Government & Civil Society"; $data[71,1] = "Conflict, Peace & Security"; $data[71,2] = "150"; $data[71,3] = "152"
$data[72,0] = "Government & Civil Society"; $data[72,1] = "Conflict, Peace & Security"; $data[72,2] = "150"; $data[72,3] = "152"
$data[73,0] = "Other Social Infrastructure & Services"; $data[73,1] = "Other Social Infrastructure & Services"; $data[73,2] = "160"; $data[73,3] = "160"
$data[74,0] = "Other Social Infrastructure & Services"; $data[74,1] = "Other Social Infrastructure & Services"; $data[74,2] = "160"; $data[74,3] = "160"
$data[75,0] = "Other Social Infrastructure & Services"; $data[75,1] = "Other Social Infrastructure & Services"; $data[75,2] = "160"; $data[75,3] = "160"
$data[76,0] = "Other Social Infrastructure & Services"; $data[76,1] = "Other Social Infrastructure & Services"; $data[76,2] = "160"; $data[76,3] = "160"
$data[77,0] = "Other Social Infrastructure & Services"; $data[77,1] = "Other Social Infrastructure & Services"; $data[77,2] = "160"; $data[77,3] = "160"
$data[78,0] = "Other Social Infrastructure & Services"; $data[78,1] = "Other Social Infrastructure & Services"; $data[78,2] = "160"; $data[78,3] = "160"
$data[79,0] = "Other Social Infrastructure & Services"; $data[79,1] = "Other Social Infrastructure & Services"; $data[79,2] = "160"; $data[79,3] = "160"
$data[80,0] = "Other Social Infrastructure & Services"; $data[80,1] = "Other Social Infrastructure & Services"; $data[80,2] = "160"; $data[80,3] = "160"
$data[81,0] = "Other Social Infrastructure & Services"; $data[81,1] = "Other Social Infrastructure & Services"; $data[81,2] = "160"; $data[81,3] = "160"
$data[82,0] = "Other Social Infrastructure & Services"; $data[82,1] = "Other Social Infrastructure & Services"; $data[82,2] = "160"; $data[82,3] = "160"
$data[83,0] = "Other Social Infrastructure & Services"; $data[83,1] = "Other Social Infrastructure & Services"; $data[83,2] = "160"; $data[83,3] = "160"
$data[84,0] = "Transport & Storage"; $data[84,1] = "Transport & Storage"; $data[84,2] = "210"; $data[84,3] = "210"
$data[85,0] = "Transport & Storage"; $data[85,1] = "Transport & Storage"; $data[85,2] = "210"; $data[85,3] = "210"
$data[86,0] = "Transport & Storage"; $data[86,1] = "Transport & Storage"; $data[86,2] = "210"; $data[86,3] = "210"
$data[87,0] = "Transport & Storage"; $data[87,1] = "Transport & Storage"; $data[87,2] = "210"; $data[87,3] = "210"
$data[88,0] = "Transport & Storage"; $data[88,1] = "Transport & Storage"; $data[88,2] = "210"; $data[88,3] = "210"
$data[89,0] = "Transport & Storage"; $data[89,1] = "Transport & Storage"; $data[89,2] = "210"; $data[89,3] = "210"
$data[90,0] = "Transport & Storage"; $data[90,1] = "Transport & Storage"; $data[90,2] = "210"; $data[90,3] = "210"
$data[91,0] = "Communications"; $data[91,1] = "Communications"; $data[91,2] = "220"; $data[91,3] = "220"
$data[92,0] = "Communications"; $data[92,1] = "Communications"; $data[92,2] = "220"; $data[92,3] = "220"
$data[93,0] = "Communications"; $data[93,1] = "Communications"; $data[93,2] = "220"; $data[93,3] = "220"
$data[94,0] = "Communications"; $data[94,1] = "Communications"; $data[94,2] = "220"; $data[94,3] = "220"
$data[95,0] = "Energy"; $data[95,1] = "Energy Policy"; $data[95,2] = "230"; $data[95,3] = "231"
$data[96,0] = "Energy"; $data[96,1] = "Energy Policy"; $data[96,2] = "230"; $data[96,3] = "231"
$data[97,0] = "Energy"; $data[97,1] = "Energy Policy"; $data[97,2] = "230"; $data[97,3] = "231"
$data[98,0] = "Energy"; $data[98,1] = "Energy Policy"; $data[98,2] = "230"; $data[98,3] = "231"
$data[99,0] = "Energy"; $data[99,1] = "Energy generation, renewable sources"; $data[99,2] = "230"; $data[99,3] = "232"
$data[100,0] = "Energy"; $data[100,1] = "Energy generation, renewable sources"; $data[100,2] = "230"; $data[100,3] = "232"
$data[101,0] = "Energy"; $data[101,1] = "Energy generation, renewable sources"; $data[101,2] = "230"; $data[101,3] = "232"
$data[102,0] = "Energy"; $data[102,1] = "Energy generation, renewable sources"; $data[102,2] = "230"; $data[102,3] = "232"
$data[103,0] = "Energy"; $data[103,1] = "Energy generation, renewable sources"; $data[103,2] = "230"; $data[103,3] = "232"
$data[104,0] = "Energy"; $data[104,1] = "Energy generation, renewable sources"; $data[104,2] = "230"; $data[104,3] = "232"
$data[105,0] = "Energy"; $data[105,1] = "Energy generation, renewable sources"; $data[105,2] = "230"; $data[105,3] = "232"
$data[106,0] = "Energy"; $data[106,1] = "Energy generation, renewable sources"; $data[106,2] = "230"; $data[106,3] = "232"
$data[107,0] = "Energy"; $data[107,1] = "Energy generation, renewable sources"; $data[107,2] = "230"; $data[107,3] = "232"
$data[108,0] = "Energy"; $data[108,1] = "Energy generation, non-renewable sources"; $data[108,2] = "230"; $data[108,3] = "233"
$data[109,0] = "Energy"; $data[109,1] = "Energy generation, non-renewable sources"; $data[109,2] = "230"; $data[109,3] = "233"
$data[110,0] = "Energy"; $data[110,1] = "Energy generation, non-renewable sources"; $data[110,2] = "230"; $data[110,3] = "233"
$data[111,0] = "Energy"; $data[111,1] = "Energy generation, non-renewable sources"; $data[111,2] = "230"; $data[111,3] = "233"
$data[112,0] = "Energy"; $data[112,1] = "Energy generation, non-renewable sources"; $data[112,2] = "230"; $data[112,3] = "233"
$data[113,0] = "Energy"; $data[113,1] = "Energy generation, non-renewable sources"; $data[113,2] = "230"; $data[113,3] = "233"
$data[114,0] = "Energy"; $data[114,1] = "Hybrid energy plants"; $data[114,2] = "230"; $data[114,3] = "234"
$data[115,0] = "Energy"; $data[115,1] = "Nuclear energy plants"; $data[115,2] = "230"; $data[115,3] = "235"
$data[116,0] = "Energy"; $data[116,1] = "Energy distribution"; $data[116,2] = "230"; $data[116,3] = "236"
$data[117,0] = "Energy"; $data[117,1] = "Energy distribution"; $data[117,2] = "230"; $data[117,3] = "236"
$data[118,0] = "Energy"; $data[118,1] = "Energy distribution"; $data[118,2] = "230"; $data[118,3] = "236"
$data[119,0] = "Energy"; $data[119,1] = "Energy distribution"; $data[119,2] = "230"; $data[119,3] = "236"
$data[120,0] = "Energy"; $data[120,1] = "Energy distribution"; $data[120,2] = "230"; $data[120,3] = "236"
$data[121,0] = "Energy"; $data[121,1] = "Energy distribution"; $data[121,2] = "230"; $data[121,3] = "236"
$data[122,0] = "Energy"; $data[122,1] = "Energy distribution"; $data[122,2] = "230"; $data[122,3] = "236"
$data[123,0] = "Banking & Financial Services"; $data[123,1] = "Banking & Financial Services"; $data[123,2] = "240"; $data[123,3] = "240"
$data[124,0] = "Banking & Financial Services"; $data[124,1] = "Banking & Financial Services"; $data[124,2] = "240"; $data[124,3] = "240"
$data[125,0] = "Banking & Financial Services"; $data[125,1] = "Banking & Financial Services"; $data[125,2] = "240"; $data[125,3] = "240"
$data[126,0] = "Banking & Financial Services"; $data[126,1] = "Banking & Financial Services"; $data[126,2] = "240"; $data[126,3] = "240"
$data[127,0] = "Banking & Financial Services"; $data[127,1] = "Banking & Financial Services"; $data[127,2] = "240"; $data[127,3] = "240"
$data[128,0] = "Banking & Financial Services"; $data[128,1] = "Banking & Financial Services"; $data[128,2] = "240"; $data[128,3] = "240"
$data[129,0] = "Business & Other Services"; $data[129,1] = "Business & Other Services"; $data[129,2] = "250"; $data[129,3] = "250"
$data[130,0] = "Business & Other Services"; $data[130,1] = "Business & Other Services"; $data[130,2] = "250"; $data[130,3] = "250"
$data[131,0] = "Business & Other Services"; $data[131,1] = "Business & Other Services"; $data[131,2] = "250"; $data[131,3] = "250"
$data[132,0] = "Business & Other Services"; $data[132,1] = "Business & Other Services"; $data[132,2] = "250"; $data[132,3] = "250"
$data[133,0] = "Agriculture, Forestry, Fishing"; $data[133,1] = "Agriculture"; $data[133,2] = "310"; $data[133,3] = "311"
$data[134,0] = "Agriculture, Forestry, Fishing"; $data[134,1] = "Agriculture"; $data[134,2] = "310"; $data[134,3] = "311"
$data[135,0] = "Agriculture, Forestry, Fishing"; $data[135,1] = "Agriculture"; $data[135,2] = "310"; $data[135,3] = "311"
$data[136,0] = "Agriculture, Forestry, Fishing"; $data[136,1] = "Agriculture"; $data[136,2] = "310"; $data[136,3] = "311"
$data[137,0] = "Agriculture, Forestry, Fishing"; $data[137,1] = "Agriculture"; $data[137,2] = "310"; $data[137,3] = "311"
$data[138,0] = "Agriculture, Forestry, Fishing"; $data[138,1] = "Agriculture"; $data[138,2] = "310"; $data[138,3] = "311"
$data[139,0] = "Agriculture, Forestry, Fishing"; $data[139,1] = "Agriculture"; $data[139,2] = "310"; $data[139,3] = "311"
$data[140,0] = "Agriculture, Forestry, Fishing"; $data[140,1] = "Agriculture"; $data[140,2] = "310"; $data[140,3] = "311"
$data[141,0] = "Agriculture, Forestry, Fishing"; $data[141,1] = "Agriculture"; $data[141,2] = "310"; $data[141,3] = "311"
$data[142,0] = "Agriculture, Forestry, Fishing"; $data[142,1] = "Agriculture"; $data[142,2] = "310"; $data[142,3] = "311"
$data[143,0] = "Agriculture, Forestry, Fishing"; $data[143,1] = "Agriculture"; $data[143,2] = "310"; $data[143,3] = "311"
$data[144,0] = "Agriculture, Forestry, Fishing"; $data[144,1] = "Agriculture"; $data[144,2] = "310"; $data[144,3] = "311"
$data[145,0] = "Agriculture, Forestry, Fishing"; $data[145,1] = "Agriculture"; $data[145,2] = "310"; $data[145,3] = "311"
$data[146,0] = "Agriculture, Forestry, Fishing"; $data[146,1] = "Agriculture"; $data[146,2] = "310"; $data[146,3] = "311"
$data[147,0] = "Agriculture, Forestry, Fishing"; $data[147,1] = "Agriculture"; $data[147,2] = "310"; $data[147,3] = "311"
$data[148,0] = "Agriculture, Forestry, Fishing"; $data[148,1] = "Agriculture"; $data[148,2] = "310"; $data[148,3] = "311"
$data[149,0] = "Agriculture, Forestry, Fishing"; $data[149,1] = "Agriculture"; $data[149,2] = "310"; $data[149,3] = "311"
$data[150,0] = "Agriculture, Forestry, Fishing"; $data[150,1] = "Agriculture"; $data[150,2] = "310"; $data[150,3] = "311"
$data[151,0] = "Agriculture, Forestry, Fishing"; $data[151,1] = "Forestry"; $data[151,2] = "310"; $data[151,3] = "312"
$data[152,0] = "Agriculture, Forestry, Fishing"; $data[152,1] = "Forestry"; $data[152,2] = "310"; $data[152,3] = "312"
$data[153,0] = "Agriculture, Forestry, Fishing"; $data[153,1] = "Forestry"; $data[153,2] = "310"; $data[153,3] = "312"
$data[154,0] = "Agriculture, Forestry, Fishing"; $data[154,1] = "Forestry"; $data[154,2] = "310"; $data[154,3] = "312"
$data[155,0] = "Agriculture, Forestry, Fishing"; $data[155,1] = "Forestry"; $data[155,2] = "310"; $data[155,3] = "312"
$data[156,0] = "Agriculture, Forestry, Fishing"; $data[156,1] = "Forestry"; $data[156,2] = "310"; $data[156,3] = "312"
$data[157,0] = "Agriculture, Forestry, Fishing"; $data[157,1] = "Fishing"; $data[157,2] = "310"; $data[157,3] = "313"
$data[158,0] = "Agriculture, Forestry, Fishing"; $data[158,1] = "Fishing"; $data[158,2] = "310"; $data[158,3] = "313"
$data[159,0] = "Agriculture, Forestry, Fishing"; $data[159,1] = "Fishing"; $data[159,2] = "310"; $data[159,3] = "313"
$data[160,0] = "Agriculture, Forestry, Fishing"; $data[160,1] = "Fishing"; $data[160,2] = "310"; $data[160,3] = "313"
$data[161,0] = "Agriculture, Forestry, Fishing"; $data[161,1] = "Fishing"; $data[161,2] = "310"; $data[161,3] = "313"
$data[162,0] = "Industry, Mining, Construction"; $data[162,1] = "Industry"; $data[162,2] = "320"; $data[162,3] = "321"
$data[163,0] = "Industry, Mining, Construction"; $data[163,1] = "Industry"; $data[163,2] = "320"; $data[163,3] = "321"
$data[164,0] = "Industry, Mining, Construction"; $data[164,1] = "Industry"; $data[164,2] = "320"; $data[164,3] = "321"
$data[165,0] = "Industry, Mining, Construction"; $data[165,1] = "Industry"; $data[165,2] = "320"; $data[165,3] = "321"
$data[166,0] = "Industry, Mining, Construction"; $data[166,1] = "Industry"; $data[166,2] = "320"; $data[166,3] = "321"
$data[167,0] = "Industry, Mining, Construction"; $data[167,1] = "Industry"; $data[167,2] = "320"; $data[167,3] = "321"
$data[168,0] = "Industry, Mining, Construction"; $data[168,1] = "Industry"; $data[168,2] = "320"; $data[168,3] = "321"
$data[169,0] = "Industry, Mining, Construction"; $data[169,1] = "Industry"; $data[169,2] = "320"; $data[169,3] = "321"
$data[170,0] = "Industry, Mining, Construction"; $data[170,1] = "Industry"; $data[170,2] = "320"; $data[170,3] = "321"
$data[171,0] = "Industry, Mining, Construction"; $data[171,1] = "Industry"; $data[171,2] = "320"; $data[171,3] = "321"
$data[172,0] = "Industry, Mining, Construction"; $data[172,1] = "Industry"; $data[172,2] = "320"; $data[172,3] = "321"
$data[173,0] = "Industry, Mining, Construction"; $data[173,1] = "Industry"; $data[173,2] = "320"; $data[173,3] = "321"
$data[174,0] = "Industry, Mining, Construction"; $data[174,1] = "Industry"; $data[174,2] = "320"; $data[174,3] = "321"
$data[175,0] = "Industry, Mining, Construction"; $data[175,1] = "Industry"; $data[175,2] = "320"; $data[175,3] = "321"
$data[176,0] = "Industry, Mining, Construction"; $data[176,1] = "Industry"; $data[176,2] = "320"; $data[176,3] = "321"
$data[177,0] = "Industry, Mining, Construction"; $data[177,1] = "Industry"; $data[177,2] = "320"; $data[177,3] = "321"
$data[178,0] = "Industry, Mining, Construction"; $data[178,1] = "Industry"; $data[178,2] = "320"; $data[178,3] = "321"
$data[179,0] = "Industry, Mining, Construction"; $data[179,1] = "Industry"; $data[179,2] = "320"; $data[179,3] = "321"
$data[180,0] = "Industry, Mining, Construction"; $data[180,1] = "Industry"; $data[180,2] = "320"; $data[180,3] = "321"
$data[181,0] = "Industry, Mining, Construction"; $data[181,1] = "Mineral Resources & Mining"; $data[181,2] = "320"; $data[181,3] = "322"
$data[182,0] = "Industry, Mining, Construction"; $data[182,1] = "Mineral Resources & Mining"; $data[182,2] = "320"; $data[182,3] = "322"
$data[183,0] = "Industry, Mining, Construction"; $data[183,1] = "Mineral Resources & Mining"; $data[183,2] = "320"; $data[183,3] = "322"
$data[184,0] = "Industry, Mining, Construction"; $data[184,1] = "Mineral Resources & Mining"; $data[184,2] = "320"; $data[184,3] = "322"
$data[185,0] = "Industry, Mining, Construction"; $data[185,1] = "Mineral Resources & Mining"; $data[185,2] = "320"; $data[185,3] = "322"
$data[186,0] = "Industry, Mining, Construction"; $data[186,1] = "Mineral Resources & Mining"; $data[186,2] = "320"; $data[186,3] = "322"
$data[187,0] = "Industry, Mining, Construction"; $data[187,1] = "Mineral Resources & Mining"; $data[187,2] = "320"; $data[187,3] = "322"
$data[188,0] = "Industry, Mining, Construction"; $data[188,1] = "Mineral Resources & Mining"; $data[188,2] = "320"; $data[188,3] = "322"
$data[189,0] = "Industry, Mining, Construction"; $data[189,1] = "Mineral Resources & Mining"; $data[189,2] = "320"; $data[189,3] = "322"
$data[190,0] = "Industry, Mining, Construction"; $data[190,1] = "Mineral Resources & Mining"; $data[190,2] = "320"; $data[190,3] = "322"
$data[191,0] = "Industry, Mining, Construction"; $data[191,1] = "Construction"; $data[191,2] = "320"; $data[191,3] = "323"
$data[192,0] = "Trade Policies & Regulations"; $data[192,1] = "Trade Policies & Regulations"; $data[192,2] = "331"; $data[192,3] = "331"
$data[193,0] = "Trade Policies & Regulations"; $data[193,1] = "Trade Policies & Regulations"; $data[193,2] = "331"; $data[193,3] = "331"
$data[194,0] = "Trade Policies & Regulations"; $data[194,1] = "Trade Policies & Regulations"; $data[194,2] = "331"; $data[194,3] = "331"
$data[195,0] = "Trade Policies & Regulations"; $data[195,1] = "Trade Policies & Regulations"; $data[195,2] = "331"; $data[195,3] = "331"
$data[196,0] = "Trade Policies & Regulations"; $data[196,1] = "Trade Policies & Regulations"; $data[196,2] = "331"; $data[196,3] = "331"
$data[197,0] = "Trade Policies & Regulations"; $data[197,1] = "Trade Policies & Regulations"; $data[197,2] = "331"; $data[197,3] = "331"
$data[198,0] = "Tourism"; $data[198,1] = "Tourism"; $data[198,2] = "332"; $data[198,3] = "332"
$data[199,0] = "General Environment Protection"; $data[199,1] = "General Environment Protection"; $data[199,2] = "410"; $data[199,3] = "410"
$data[200,0] = "General Environment Protection"; $data[200,1] = "General Environment Protection"; $data[200,2] = "410"; $data[200,3] = "410"
$data[201,0] = "General Environment Protection"; $data[201,1] = "General Environment Protection"; $data[201,2] = "410"; $data[201,3] = "410"
$data[202,0] = "General Environment Protection"; $data[202,1] = "General Environment Protection"; $data[202,2] = "410"; $data[202,3] = "410"
$data[203,0] = "General Environment Protection"; $data[203,1] = "General Environment Protection"; $data[203,2] = "410"; $data[203,3] = "410"
$data[204,0] = "General Environment Protection"; $data[204,1] = "General Environment Protection"; $data[204,2] = "410"; $data[204,3] = "410"
$data[205,0] = "Other Multisector"; $data[205,1] = "Other Multisector"; $data[205,2] = "430"; $data[205,3] = "430"
$data[206,0] = "Other Multisector"; $data[206,1] = "Other Multisector"; $data[206,2] = "430"; $data[206,3] = "430"
$data[207,0] = "Other Multisector"; $data[207,1] = "Other Multisector"; $data[207,2] = "430"; $data[207,3] = "430"
$data[208,0] = "Other Multisector"; $data[208,1] = "Other Multisector"; $data[208,2] = "430"; $data[208,3] = "430"
$data[209,0] = "Other Multisector"; $data[209,1] = "Other Multisector"; $data[209,2] = "430"; $data[209,3] = "430"
$data[210,0] = "Other Multisector"; $data[210,1] = "Other Multisector"; $data[210,2] = "430"; $data[210,3] = "430"
$data[211,0] = "Other Multisector"; $data[211,1] = "Other Multisector"; $data[211,2] = "430"; $data[211,3] = "430"
$data[212,0] = "Other Multisector"; $data[212,1] = "Other Multisector"; $data[212,2] = "430"; $data[212,3] = "430"
$data[213,0] = "Other Multisector"; $data[213,1] = "Other Multisector"; $data[213,2] = "430"; $data[213,3] = "430"
$data[214,0] = "Other Multisector"; $data[214,1] = "Other Multisector"; $data[214,2] = "430"; $data[214,3] = "430"
$data[215,0] = "General Budget Support"; $data[215,1] = "General Budget Support"; $data[215,2] = "510"; $data[215,3] = "510"
$data[216,0] = "Development Food Assistance"; $data[216,1] = "Development Food Assistance"; $data[216,2] = "520"; $data[216,3] = "520"
$data[217,0] = "Other Commodity Assistance"; $data[217,1] = "Other Commodity Assistance"; $data[217,2] = "530"; $data[217,3] = "530"
$data[218,0] = "Other Commodity Assistance"; $data[218,1] = "Other Commodity Assistance"; $data[218,2] = "530"; $data[218,3] = "530"
$data[219,0] = "Action Relating to Debt"; $data[219,1] = "Action Relating to Debt"; $data[219,2] = "600"; $data[219,3] = "600"
$data[220,0] = "Action Relating to Debt"; $data[220,1] = "Action Relating to Debt"; $data[220,2] = "600"; $data[220,3] = "600"
$data[221,0] = "Action Relating to Debt"; $data[221,1] = "Action Relating to Debt"; $data[221,2] = "600"; $data[221,3] = "600"
$data[222,0] = "Action Relating to Debt"; $data[222,1] = "Action Relating to Debt"; $data[222,2] = "600"; $data[222,3] = "600"
$data[223,0] = "Action Relating to Debt"; $data[223,1] = "Action Relating to Debt"; $data[223,2] = "600"; $data[223,3] = "600"
$data[224,0] = "Action Relating to Debt"; $data[224,1] = "Action Relating to Debt"; $data[224,2] = "600"; $data[224,3] = "600"
$data[225,0] = "Action Relating to Debt"; $data[225,1] = "Action Relating to Debt"; $data[225,2] = "600"; $data[225,3] = "600"
$data[226,0] = "Emergency Response"; $data[226,1] = "Emergency Response"; $data[226,2] = "720"; $data[226,3] = "720"
$data[227,0] = "Emergency Response"; $data[227,1] = "Emergency Response"; $data[227,2] = "720"; $data[227,3] = "720"
$data[228,0] = "Emergency Response"; $data[228,1] = "Emergency Response"; $data[228,2] = "720"; $data[228,3] = "720"
$data[229,0] = "Reconstruction Relief & Rehabilitation"; $data[229,1] = "Reconstruction Relief & Rehabilitation"; $data[229,2] = "730"; $data[229,3] = "730"
$data[230,0] = "Disaster Prevention & Preparedness"; $data[230,1] = "Disaster Prevention & Preparedness"; $data[230,2] = "740"; $data[230,3] = "740"
$data[231,0] = "Administrative Costs of Donors"; $data[231,1] = "Administrative Costs of Donors"; $data[231,2] = "910"; $data[231,3] = "910"
$data[232,0] = "Refugees in Donor Countries"; $data[232,1] = "Refugees in Donor Countries"; $data[232,2] = "930"; $data[232,3] = "930"
$data[233,0] = "Unallocated / Unspecified"; $data[233,1] = "Unallocated / Unspecified"; $data[233,2] = "998"; $data[233,3] = "998"
$data[234,0] = "Unallocated / Unspecified"; $data[234,1] = "Unallocated / Unspecified"; $data[234,2] = "998"; $data[234,3] = "998"

$ws.Range("D1:G235").Value = $data
